# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect freshly generated data (gh-pages output at commit a3196b5).

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 2128
    3 = 617
    4 = 1511
    5 = 7184
    6 = 176
    7 = 144
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
